# Protocol.docx edit script
# Implements the changes described by the commit:
# "Refitting data with the new algo that includes the Ac + Ap case."

$word.UserName = "Julien LAMOUR"
$word.UserInitials = "JL"

$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. "27 species were selected " -> "Twenty-seven species were selected "
# ------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("27 species were selected", $false, $false, $false, $false, $false, $true, 1, $false, "Twenty-seven species were selected", 2) | Out-Null

# ------------------------------------------------------------------
# 2. Remove the old "_GoBack" bookmark that sits around "in-situ " —
#    it will be superseded later (Word relocates it to the very last
#    edit location, which ends up inside the new comment below).
# ------------------------------------------------------------------
try {
    $bm = $d.Bookmarks("_GoBack")
    $bm.Delete()
} catch {
}

# ------------------------------------------------------------------
# 3. Drop the stray " (Figure 1)" after "8 am to 4 pm" and re-insert it
#    after "and bending the stem" instead.
# ------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute(" Measurements were performed from 8 am to 4 pm (Figure 1). Each leaf", $false, $false, $false, $false, $false, $true, 1, $false, " Measurements were performed from 8 am to 4 pm. Each leaf", 2) | Out-Null

$rng = $d.Content
$rng.Find.Execute(" and bending the stem. All the leaves", $false, $false, $false, $false, $false, $true, 1, $false, " and bending the stem (Figure 1). All the leaves", 2) | Out-Null

# ------------------------------------------------------------------
# 4. Remove the bookmark-spanning run split around "in-situ (n = 65)"
#    (merge the trailing space into the following run).
# ------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Most gas exchange measurements were performed in-situ (n = 65)", $false, $false, $false, $false, $false, $true, 1, $false, "Most gas exchange measurements were performed in-situ (n = 65)", 2) | Out-Null

# ------------------------------------------------------------------
# 5. Insert "corrected " before "spectra per leaf were averaged"
# ------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("The 3 spectra per leaf were averaged", $false, $false, $false, $false, $false, $true, 1, $false, "The 3 corrected spectra per leaf were averaged", 2) | Out-Null

# ------------------------------------------------------------------
# 6. Add the two reviewer comments.
# ------------------------------------------------------------------

# Comment 0: around "typical old-growth forest species (e.g. Bocoa prouacensis)"
$rng = $d.Content
$rng.Find.Execute("typical old-growth forest species (e.g. Bocoa prouacensis)", $false) | Out-Null
$comment0 = $d.Comments.Add($rng, "Est ce que vous pouvez préciser les strategies des espèces si vous les connaissez?")

# Comment 1: around the full "Leaf elemental composition ... Germany)." sentence
$rng = $d.Content
$rng.Find.Execute("Leaf elemental composition was measured by the Silva lab (Nancy, France) with the Unicube elemental analyzer (Elementar, Langenselbold, Germany).", $false) | Out-Null
$comment1 = $d.Comments.Add($rng, "True?")

# ------------------------------------------------------------------
# 7. Best-effort: materialize the (otherwise dangling) comment styles
#    referenced by the comment runs above.
# ------------------------------------------------------------------
try { $d.Styles.Add("Comment Reference", 2) | Out-Null } catch {}
try { $d.Styles.Add("Comment Text", 1) | Out-Null } catch {}

Write-Output "done"
